$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "mgrael"
$ws.Range("B7").Value = "enkwudhu"
$ws.Range("C7").Value = "rmgraelenkwudhu@gmail.com"
$ws.Range("D7").Value = "YTWYRlMctmNuilk"
$ws.Range("E7").Value = "17-02-2022 08:13"

$ws.Range("A8").Value = "fzydmi"
$ws.Range("B8").Value = "hieekebs"
$ws.Range("C8").Value = "yfzydmihieekebs@gmail.com"
$ws.Range("D8").Value = "fnxeFuwmIyFxwAf"
$ws.Range("E8").Value = "17-02-2022 08:14"
